$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (A1:C1 relabeled, D1 added) - exporting with a Schema now
# adds a 4th "Ciudad" column.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Edad"
$ws.Range("C1").Value = "Telefono"
$ws.Range("D1").Value = "Ciudad"

# Header row is shorter now (was 18pt, now 15pt).
$ws.Rows.Item(1).RowHeight = 15.0

# Columns 1, 2 and the new column 4 go back to the sheet's default width;
# column 3 keeps a slightly wider, but still near-default, custom width.
$ws.Columns.Item(1).ColumnWidth = 8.333333333333334
$ws.Columns.Item(2).ColumnWidth = 8.333333333333334
$ws.Columns.Item(3).ColumnWidth = 8.833333333333334
$ws.Columns.Item(4).ColumnWidth = 8.333333333333334
